$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for columns A and B
$ws.Range("A:B").ColumnWidth = 16.42578125

# Update cell values (rows 1-4), row 5 remains unchanged
$ws.Range("A1").Value = -0.021519257133428595
$ws.Range("B1").Value = -0.021122500817468581
$ws.Range("A2").Value = -0.014013515714812306
$ws.Range("B2").Value = -0.040831184676669942
$ws.Range("A3").Value = -0.00072270596763810313
$ws.Range("B3").Value = -0.00057434853241421148
$ws.Range("A4").Value = -0.088976408078623087
$ws.Range("B4").Value = -0.088894454114066823
